$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "231×6=1386" "972×6=5832"
Replace-Text "392×9=3528" "902×8=7216"
Replace-Text "102×2=204" "548×9=4932"
Replace-Text "573×5=2865" "648×4=2592"
Replace-Text "762×9=6858" "468×4=1872"
Replace-Text "342×6=2052" "740×6=4440"
Replace-Text "759×8=6072" "133×3=399"
Replace-Text "490×3=1470" "844×9=7596"
Replace-Text "547×6=3282" "495×4=1980"
Replace-Text "368×3=1104" "180×8=1440"
Replace-Text "816×7=5712" "716×2=1432"
Replace-Text "450×5=2250" "877×7=6139"
Replace-Text "999×2=1998" "231×4=924"
Replace-Text "917×9=8253" "975×9=8775"
Replace-Text "512×4=2048" "296×7=2072"
Replace-Text "268×9=2412" "897×7=6279"
Replace-Text "266×7=1862" "452×8=3616"
Replace-Text "612×4=2448" "785×7=5495"
Replace-Text "433×7=3031" "345×8=2760"
Replace-Text "616×6=3696" "564×7=3948"
Replace-Text "808×7=5656" "236×3=708"
Replace-Text "765×7=5355" "892×5=4460"
Replace-Text "530×7=3710" "171×6=1026"
Replace-Text "177×4=708" "772×7=5404"
Replace-Text "431×7=3017" "370×7=2590"
